# "19/04 -> Doing the same with ProductionActivity"
#
# - fills in the empty paragraph right after the "18/04" fragment
#   discussion with "Create CellObject class"
# - appends five new paragraphs describing the 19/04 CellAdapter
#   colour work, a bulleted "Doing the same with ProductionActivity"
#   note (carrying the relocated "_GoBack" bookmark) and a final
#   "Rafraichir" link note.
$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "Create CellObject class" in the previously empty paragraph.
# ------------------------------------------------------------------
$pCreateCell = $d.Paragraphs(15)
$pCreateCell.Range.Text = "Create CellObject class"

# ------------------------------------------------------------------
# 2. Mint the (localised) "List Paragraph" style definition into
#    styles.xml via a throwaway paragraph, without disturbing the
#    real document content.
# ------------------------------------------------------------------
$mintAnchor = $d.Paragraphs($d.Paragraphs.Count)
$mintRange = $mintAnchor.Range
$mintRange.Collapse(0)
$mintRange.InsertParagraphAfter()
$mintPara = $d.Paragraphs($d.Paragraphs.Count)
$mintPara.Style = "List Paragraph"
$mintPara.Style = "Normal"
$mintPara.Range.Delete()

# ------------------------------------------------------------------
# 3. Append the five new tail paragraphs after "Add all the colour
#    corresponding to color code in the cellAdapter."
# ------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$pLast = $d.Paragraphs($lastIndex)
$tailRange = $pLast.Range
$tailRange.Collapse(0)
$tailRange.InsertParagraphAfter()
$tailRange.InsertParagraphAfter()
$tailRange.InsertParagraphAfter()
$tailRange.InsertParagraphAfter()
$tailRange.InsertParagraphAfter()

$pGalere = $d.Paragraphs($lastIndex + 1)
$pGalere.Range.InsertAfter("Beaucoup galéré : pb -> je voulais parcourir la ListView dans l’activité pour mettre à chaque fois le fond TextView du code couleur avec la couleur correspondant au code. C’était impossible, la solution a été de le faire directement dans CellAdapter au moment du remplissage des vues.")

# $lastIndex + 2 stays a blank "Courier New" paragraph.

$pProduction = $d.Paragraphs($lastIndex + 3)
$pProduction.Range.InsertAfter("Doing the same with ProductionActivity")
$pProduction.Range.ListFormat.ApplyBulletDefault()

# $lastIndex + 4 stays a blank "Courier New" paragraph.

$pRafraichir = $d.Paragraphs($lastIndex + 5)
$rRafraichir = $pRafraichir.Range
$rRafraichir.InsertAfter("Rafraichir : ")
$rRafraichir.Collapse(0)
$rRafraichir.InsertAfter("http://www.softwarepassion.com/android-series-custom-listview-items-and-adapters/")

# ------------------------------------------------------------------
# 4. Move the "_GoBack" bookmark from the end of "Add all the colour
#    ..." to the start of the new "Doing the same with
#    ProductionActivity" paragraph.
# ------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
$bmRange = $pProduction.Range
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)
